# Update roteiro de entregas (delivery route) sheet:
# - Rows 2-16: replace with the new Nov 21 route assignments
# - Rows 17-26: append the new Nov 22 deliveries
# - Column E (Bairro) widened from 13 to 15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("11081", "21/11/2024", "manhã", "Desconhecido", "Brejarú"),
    @("11081", "21/11/2024", "manhã", "Desconhecido", "Brejarú"),
    @("11191", "21/11/2024", "manhã", "Desconhecido", "Pedra Branca"),
    @("11191", "21/11/2024", "manhã", "Desconhecido", "Pedra Branca"),
    @("11192", "21/11/2024", "manhã", "Desconhecido", "Pedra Branca"),
    @("11192", "21/11/2024", "manhã", "Desconhecido", "Pedra Branca"),
    @("11190", "21/11/2024", "manhã", "Desconhecido", "Pedra Branca"),
    @("11184", "21/11/2024", "tarde", "Desconhecido", "Pagani"),
    @("11184", "21/11/2024", "tarde", "Desconhecido", "Pagani"),
    @("11190", "21/11/2024", "tarde", "Desconhecido", "Pedra Branca"),
    @("11195", "21/11/2024", "tarde", "Desconhecido", "Pedra Branca"),
    @("11195", "21/11/2024", "tarde", "Desconhecido", "Pedra Branca"),
    @("11189", "21/11/2024", "tarde", "Desconhecido", "Pedra Branca"),
    @("11189", "21/11/2024", "tarde", "Desconhecido", "Pedra Branca"),
    @("11220", "21/11/2024", "tarde", "Desconhecido", "São Sebastião"),
    @("11141", "22/11/2024", "manhã", "Desconhecido", "Centro"),
    @("11141", "22/11/2024", "manhã", "Desconhecido", "Centro"),
    @("11145", "22/11/2024", "manhã", "Desconhecido", "Centro"),
    @("11145", "22/11/2024", "manhã", "Desconhecido", "Centro"),
    @("11214", "22/11/2024", "manhã", "Desconhecido", "Rio Grande"),
    @("11057", "22/11/2024", "tarde", "Desconhecido", "Aririú"),
    @("11057", "22/11/2024", "tarde", "Desconhecido", "Aririú"),
    @("11214", "22/11/2024", "tarde", "Desconhecido", "Rio Grande"),
    @("11208", "22/11/2024", "tarde", "Desconhecido", "Praia de Fora"),
    @("11208", "22/11/2024", "tarde", "Desconhecido", "Praia de Fora")
)

# A "template" cell that already carries the sheet's standard data-row style
# (centered, bordered-free body style used by rows 2-16) so every new/edited
# cell keeps that exact formatting.
$templateCell = $ws.Range("A2")

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i

    # Column A ("Nº Pedido") looks like a plain integer, so Excel would
    # normally coerce it to a Number on assignment under General formatting.
    # Force the cell to Text first, write the value, then restore the
    # original "style 3" look (General, centered) via a formats-only paste
    # from the template cell so the literal text is preserved without a
    # lingering custom number format.
    $aCell = $ws.Cells.Item($r, 1)
    $aCell.NumberFormat = "@"
    $aCell.Value = $row[0]
    $templateCell.Copy()
    $aCell.PasteSpecial(-4122)

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Widen column E (Bairro) from 13 to 15 characters.
$ws.Range("E1").ColumnWidth = 14.166666666666666

$excel.CutCopyMode = 0
